$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Dylan block (row 7): task #71 was finished as a *non*-development task,
#    re-labelled "Non development Task #88" and highlighted yellow.
# ---------------------------------------------------------------------------
$ws.Range("E7").Value = "Non development Task #88"
$ws.Range("E7").Interior.Color = 65535   # RGB(255,255,0) -> yellow

# ---------------------------------------------------------------------------
# 2) Matis block (row 12): task #72 likewise reclassified as non-development,
#    re-labelled "Non development Task #93" and recoloured to match the
#    existing "Non development" gold fill used elsewhere in the sheet.
# ---------------------------------------------------------------------------
$ws.Range("E12").Value = "Non development Task #93"
$ws.Range("E12").Interior.Color = 49407  # RGB(255,192,0) -> gold

# ---------------------------------------------------------------------------
# 3) Dylan block (row 8): the empty placeholder cell E8 is cleared completely
#    (value + formatting), shrinking that block by one visible slot.
# ---------------------------------------------------------------------------
$ws.Range("E8").Clear()

# ---------------------------------------------------------------------------
# 4) PA block: insert two new task rows (19 & 20) before the block's closing
#    border, so the block grows from rows 15-19 to rows 15-20.
#    First, lift the "closing" bottom border off what is currently the last
#    row (18) so it becomes a normal interior row ...
# ---------------------------------------------------------------------------
$ws.Range("C18").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
$ws.Range("D18").Borders.Item(9).LineStyle = -4142
$ws.Range("E18").Borders.Item(9).LineStyle = -4142

# ... then populate the new interior row 19 by cloning the formatting of an
# existing interior row (8) and filling in the new task ...
$ws.Range("C8").Copy()
$ws.Range("C19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D8").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E6").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Value = "Development Task #89"

# ... and finally add the new closing row 20, cloning the "last row of block"
# formatting (row 9, the end of the Dylan block) and giving it the new task.
$ws.Range("C9").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E12").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "Non development Task #94"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Keep the sheet's selection in sync with the new last cell, matching where
# Excel would naturally leave the cursor after typing the last new entry.
# ---------------------------------------------------------------------------
$ws.Range("E20").Select()
